# Apply trade #15 update to the live trading results workbook.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet: Summary
# ----------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.34   # Current Capital
$summary.Range("B4").Value = 0.34      # Total P&L $
$summary.Range("B5").Value = 0.45      # Total P&L %
$summary.Range("B6").Value = 15        # Total Trades
$summary.Range("B7").Value = 9         # Winning Trades
$summary.Range("B9").Value = 60        # Win Rate %

# ----------------------------------------------------------------------
# Sheet: Strategy Status (row 6 = MarketMaking)
# ----------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.34
$status.Range("D6").Value = 15
$status.Range("E6").Value = 0.34
$status.Range("F6").Value = 0.34
$status.Range("G6").Value = 60

# ----------------------------------------------------------------------
# New trade row (#15) appended as row 16 on both "All Trades" and
# "MarketMaking" sheets.
# ----------------------------------------------------------------------
$newRow = @(15, "2026-02-17", "23:53:40", "MarketMaking", "UP", 0.91, 0.96, "CLOSED", 5.4945, 0.05, 100.34, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

# Columns containing text that look like dates/times must be forced to
# Text format first so Excel does not auto-convert them into date/time
# serial numbers.
$textColumns = @(2, 3, 4, 5, 8, 15, 16)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item(16, $col)
        if ($textColumns -contains $col) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $newRow[$i]
    }
}
